# Update countries & provincias Spain
# - Re-sort Zambia above Guyana/Macao (table is kept sorted desc by "Casos totales")
#   and refresh Zambia's daily stats.
# - Refresh a handful of other countries' daily stats (rows 20, 60, 80, 115, 130).
# - Bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20 (Austria) ---
$ws.Range("B20").Value = 14286
$ws.Range("C20").Value = 60
$ws.Range("E20").Value = 5795

# --- Row 60 (Moldavia) ---
$ws.Range("E60").Value = 1720
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 43

# --- Row 80 (Eslovaquia) ---
$ws.Range("D80").Value = 151
$ws.Range("E80").Value = 710

# --- Row 115 (Isla de Man) ---
$ws.Range("B115").Value = 256
$ws.Range("C115").Value = 2
$ws.Range("D115").Value = 151
$ws.Range("E115").Value = 101
$ws.Range("G115").Value = 2
$ws.Range("H115").Value = 4

# --- Row 130 (Gibraltar) ---
$ws.Range("D130").Value = 104
$ws.Range("E130").Value = 25

# --- Rows 154-156: Zambia moves above Guyana/Macao (table sorted desc by total cases) ---
# Row 154 becomes Zambia, with refreshed stats
$ws.Range("A154").Value = "Zambia"
$ws.Range("B154").Value = 48
$ws.Range("C154").Value = 3
$ws.Range("D154").Value = 30
$ws.Range("E154").Value = 16
$ws.Range("F154").Value = 1
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 2

# Row 155 becomes Guyana, with its previous (unchanged) stats
$ws.Range("A155").Value = "Guyana"
$ws.Range("B155").Value = 47
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 8
$ws.Range("E155").Value = 33
$ws.Range("F155").Value = 5
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 6

# Row 156 becomes Macao, with its previous (unchanged) stats
$ws.Range("A156").Value = "Macao"
$ws.Range("B156").Value = 45
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 10
$ws.Range("E156").Value = 35
$ws.Range("F156").Value = 1
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 0

# --- Last updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 12:52"
